$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.841.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.472.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.471.14"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.76"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.058.71"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.471.39"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.882.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.03"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.86%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.609.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.81"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.73%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -10.32%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.24%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.77%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.462.25"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.39%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.57%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.13"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -13.08%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.39"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.896"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.94"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.30"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.23"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.98%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.32%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -10.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.966"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.02%  "
